$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as text (t="inlineStr") in the source workbook.
# For values that look like plain numbers, force the cell to Text format first so
# Excel does not silently convert the assigned string into a numeric value.

$ws.Range("D2").Value = '62.945.05'
$ws.Range("D3").Value = '3.154.09'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.44'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.24'
$ws.Range("D8").Value = '3.154.29'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.28'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000244'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.05'
$ws.Range("D15").Value = '3.672.32'
$ws.Range("D17").Value = '3.152.15'
$ws.Range("D18").Value = '62.934.73'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.83'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.92'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.700'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.71'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.96'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.01'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.93'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.52'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.77'
$ws.Range("D38").Value = '0.0₃0700'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0386'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '416.71'
$ws.Range("D42").Value = '2.954.96'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.258'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.39'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.82'

$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  -3.33%  '
$ws.Range("E14").Value = '  -2.55%  '
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("E24").Value = '  -2.51%  '
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  -5.02%  '
$ws.Range("E34").Value = '  -5.22%  '
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("E38").Value = '  -4.81%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("E40").Value = '  -4.39%  '
$ws.Range("E41").Value = '  -6.16%  '
$ws.Range("E42").Value = '  +3.12%  '
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("E44").Value = '  -7.18%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E47").Value = '  -3.14%  '
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("E50").Value = '  -5.84%  '
$ws.Range("E51").Value = '  -1.44%  '

$wb.Save()
